# Update the game data from "0022000098" (Rockets @ Mavericks) to
# "0022000589" (Wizards @ Bucks), including all derived team info and
# four-factor stats, per the commit: Resolved date issue '0#' --> '#'
#
# The GAME_ID is a zero-padded numeric string, so it must be forced to
# Text before assignment (otherwise Excel auto-converts it to a Number
# and drops the leading zero -- the very "date issue" the commit fixes).
# The NumberFormat is reset back to the default ("Normal" style) right
# after so the cell keeps its original (unformatted) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Washington Wizards ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0022000589"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = 1610612764
$ws.Range("D2").Value = "Wizards"
$ws.Range("E2").Value = "WAS"
$ws.Range("F2").Value = "Washington"

$ws.Range("H2").Value = 0.543
$ws.Range("I2").Value = 0.234
$ws.Range("J2").Value = 0.133
$ws.Range("K2").Value = 0.113
$ws.Range("L2").Value = 0.52
$ws.Range("M2").Value = 0.253
$ws.Range("N2").Value = 0.157
$ws.Range("O2").Value = 0.286

# --- Row 3: Milwaukee Bucks ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0022000589"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = 1610612749
$ws.Range("D3").Value = "Bucks"
$ws.Range("E3").Value = "MIL"
$ws.Range("F3").Value = "Milwaukee"

$ws.Range("H3").Value = 0.52
$ws.Range("I3").Value = 0.253
$ws.Range("J3").Value = 0.157
$ws.Range("K3").Value = 0.232
$ws.Range("L3").Value = 0.543
$ws.Range("M3").Value = 0.234
$ws.Range("N3").Value = 0.133
$ws.Range("O3").Value = 0.151
